$wb = $excel.ActiveWorkbook

# --- Sheet "Info" ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 5726785590522.906
$wsInfo.Range("B2").Value = 2.125999927520752

# --- Sheet "Activados" ---
$wsAct = $wb.Worksheets.Item("Activados")
for ($r = 2; $r -le 20; $r++) {
    $wsAct.Cells.Item($r, 1).Value = 1
    $wsAct.Cells.Item($r, 2).Value = ($r - 2) * 20
}

# --- Sheet "Operando" ---
$wsOp = $wb.Worksheets.Item("Operando")
for ($r = 2; $r -le 366; $r++) {
    $wsOp.Cells.Item($r, 1).Value = 1
}

# --- Sheet "Contaminantes" ---
$wsCont = $wb.Worksheets.Item("Contaminantes")
$wsCont.Range("B2").Value = 4851449043839.999
$wsCont.Range("C2").Value = 179.928
$wsCont.Range("B3").Value = 280418111999.9999
$wsCont.Range("C3").Value = 10.4
$wsCont.Range("B4").Value = 161348267519.9999
$wsCont.Range("C4").Value = 5.983999999999997
$wsCont.Range("B5").Value = 624762.9252863999
$wsCont.Range("C5").Value = 0.00002317088
$wsCont.Range("B6").Value = 433569542400
$wsCont.Range("C6").Value = 16.08
